# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the "全部类型" sheet, which duplicate the same exhibition rows.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# 展览 sheet (rows keyed by their row number)
$ws1.Range("F2").Value = 4669
$ws1.Range("F3").Value = 2720
$ws1.Range("F5").Value = 2731
$ws1.Range("F10").Value = 741
$ws1.Range("F12").Value = 205
$ws1.Range("F18").Value = 522
$ws1.Range("F19").Value = 522
$ws1.Range("F31").Value = 1417
$ws1.Range("F32").Value = 2276
$ws1.Range("F33").Value = 376
$ws1.Range("F39").Value = 767

# 全部类型 sheet (same underlying rows, different row numbers)
$ws4.Range("F2").Value = 4669
$ws4.Range("F3").Value = 2720
$ws4.Range("F4").Value = 2731
$ws4.Range("F8").Value = 741
$ws4.Range("F10").Value = 205
$ws4.Range("F16").Value = 522
$ws4.Range("F17").Value = 522
$ws4.Range("F33").Value = 2276
$ws4.Range("F34").Value = 376
$ws4.Range("F43").Value = 767
